# Update ligand/receptor TPM-derived expression stats and downstream
# specificity / edge-weight columns (G:J, M:T) on Sheet1 for the
# Adam9 -> Itga6 ligand-receptor pair, reflecting the refreshed TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 9.462749333333333
$ws.Range("H2").Value = 28.388248
$ws.Range("I2").Value = 0.07254428564686972
$ws.Range("J2").Value = 0.07439525120506714
$ws.Range("M2").Value = 145.7087706666667
$ws.Range("N2").Value = 437.126312
$ws.Range("O2").Value = 0.5445232453600627
$ws.Range("P2").Value = 0.5461141113270247
$ws.Range("Q2").Value = 1378.805572486819
$ws.Range("R2").Value = 12409.25015238137
$ws.Range("S2").Value = 0.03950204985276091
$ws.Range("T2").Value = 0.040628296498806

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 9.462749333333333
$ws.Range("H3").Value = 28.388248
$ws.Range("I3").Value = 0.07254428564686972
$ws.Range("J3").Value = 0.07439525120506714
$ws.Range("O3").Value = 0.003349722554576428
$ws.Range("P3").Value = 0.003359509023117945
$ws.Range("Q3").Value = 8.481944827682666
$ws.Range("R3").Value = 76.337503449144
$ws.Range("S3").Value = 0.0002430032298369546
$ws.Range("T3").Value = 0.0002499315177005492

# Row 4: ECs -> Inflammatory-Mac
$ws.Range("G4").Value = 9.462749333333333
$ws.Range("H4").Value = 28.388248
$ws.Range("I4").Value = 0.07254428564686972
$ws.Range("J4").Value = 0.07439525120506714
$ws.Range("M4").Value = 51.59199533333333
$ws.Range("N4").Value = 154.775986
$ws.Range("O4").Value = 0.1928026748491032
$ws.Range("P4").Value = 0.1933659624890163
$ws.Range("Q4").Value = 488.2021194458364
$ws.Range("R4").Value = 4393.819075012527
$ws.Range("S4").Value = 0.01398673231773389
$ws.Range("T4").Value = 0.01438550935387996

# Row 5: ECs -> MuSCs
$ws.Range("G5").Value = 9.462749333333333
$ws.Range("H5").Value = 28.388248
$ws.Range("I5").Value = 0.07254428564686972
$ws.Range("J5").Value = 0.07439525120506714
$ws.Range("M5").Value = 2.338518
$ws.Range("N5").Value = 4.677036
$ws.Range("O5").Value = 0.0087391953474509
$ws.Range("P5").Value = 0.005843151713055659
$ws.Range("Q5").Value = 22.128809645488
$ws.Range("R5").Value = 132.772857872928
$ws.Range("S5").Value = 0.000633978683609273
$ws.Range("T5").Value = 0.0004347027395220941

# Row 6: ECs -> Resolving-Mac
$ws.Range("G6").Value = 9.462749333333333
$ws.Range("H6").Value = 28.388248
$ws.Range("I6").Value = 0.07254428564686972
$ws.Range("J6").Value = 0.07439525120506714
$ws.Range("M6").Value = 67.05398933333333
$ws.Range("N6").Value = 201.161968
$ws.Range("O6").Value = 0.2505851618888069
$ws.Range("P6").Value = 0.2513172654477853
$ws.Range("Q6").Value = 634.5150928613405
$ws.Range("R6").Value = 5710.635835752064
$ws.Range("S6").Value = 0.0181785215629287
$ws.Range("T6").Value = 0.01869681109515852

# Row 7: FAPs -> ECs
$ws.Range("I7").Value = 0.3231336970688258
$ws.Range("J7").Value = 0.3313784449305509
$ws.Range("M7").Value = 145.7087706666667
$ws.Range("N7").Value = 437.126312
$ws.Range("O7").Value = 0.5445232453600627
$ws.Range("P7").Value = 0.5461141113270247
$ws.Range("Q7").Value = 6141.607684243422
$ws.Range("R7").Value = 55274.4691581908
$ws.Range("S7").Value = 0.1759538094131124
$ws.Range("T7").Value = 0.1809704449661792

# Row 8: FAPs -> FAPs
$ws.Range("I8").Value = 0.3231336970688258
$ws.Range("J8").Value = 0.3313784449305509
$ws.Range("O8").Value = 0.003349722554576428
$ws.Range("P8").Value = 0.003359509023117945
$ws.Range("S8").Value = 0.001082408233215113
$ws.Range("T8").Value = 0.001113268875810979

# Row 9: FAPs -> Inflammatory-Mac
$ws.Range("I9").Value = 0.3231336970688258
$ws.Range("J9").Value = 0.3313784449305509
$ws.Range("M9").Value = 51.59199533333333
$ws.Range("N9").Value = 154.775986
$ws.Range("O9").Value = 0.1928026748491032
$ws.Range("P9").Value = 0.1933659624890163
$ws.Range("Q9").Value = 2174.596584233878
$ws.Range("R9").Value = 19571.3692581049
$ws.Range("S9").Value = 0.06230104112874944
$ws.Range("T9").Value = 0.06407731195210944

# Row 10: FAPs -> MuSCs
$ws.Range("I10").Value = 0.3231336970688258
$ws.Range("J10").Value = 0.3313784449305509
$ws.Range("M10").Value = 2.338518
$ws.Range("N10").Value = 4.677036
$ws.Range("O10").Value = 0.0087391953474509
$ws.Range("P10").Value = 0.005843151713055659
$ws.Range("Q10").Value = 98.56826087290001
$ws.Range("R10").Value = 591.4095652374001
$ws.Range("S10").Value = 0.002823928502028491
$ws.Range("T10").Value = 0.001936294528165669

# Row 11: FAPs -> Resolving-Mac
$ws.Range("I11").Value = 0.3231336970688258
$ws.Range("J11").Value = 0.3313784449305509
$ws.Range("M11").Value = 67.05398933333333
$ws.Range("N11").Value = 201.161968
$ws.Range("O11").Value = 0.2505851618888069
$ws.Range("P11").Value = 0.2513172654477853
$ws.Range("Q11").Value = 2826.317827434578
$ws.Range("R11").Value = 25436.8604469112
$ws.Range("S11").Value = 0.0809725097917204
$ws.Range("T11").Value = 0.08328112460828554

# Row 12: Inflammatory-Mac -> ECs
$ws.Range("G12").Value = 32.300192
$ws.Range("H12").Value = 96.900576
$ws.Range("I12").Value = 0.2476229975407503
$ws.Range("J12").Value = 0.2539410918713864
$ws.Range("M12").Value = 145.7087706666667
$ws.Range("N12").Value = 437.126312
$ws.Range("O12").Value = 0.5445232453600627
$ws.Range("P12").Value = 0.5461141113270247
$ws.Range("Q12").Value = 4706.421268617301
$ws.Range("R12").Value = 42357.79141755571
$ws.Range("S12").Value = 0.1348364782466762
$ws.Range("T12").Value = 0.1386808137167565

# Row 13: Inflammatory-Mac -> FAPs
$ws.Range("G13").Value = 32.300192
$ws.Range("H13").Value = 96.900576
$ws.Range("I13").Value = 0.2476229975407503
$ws.Range("J13").Value = 0.2539410918713864
$ws.Range("O13").Value = 0.003349722554576428
$ws.Range("P13").Value = 0.003359509023117945
$ws.Range("Q13").Value = 28.952309399392
$ws.Range("R13").Value = 260.570784594528
$ws.Range("S13").Value = 0.0008294683398940747
$ws.Range("T13").Value = 0.0008531173894823456

# Row 14: Inflammatory-Mac -> Inflammatory-Mac
$ws.Range("G14").Value = 32.300192
$ws.Range("H14").Value = 96.900576
$ws.Range("I14").Value = 0.2476229975407503
$ws.Range("J14").Value = 0.2539410918713864
$ws.Range("M14").Value = 51.59199533333333
$ws.Range("N14").Value = 154.775986
$ws.Range("O14").Value = 0.1928026748491032
$ws.Range("P14").Value = 0.1933659624890163
$ws.Range("Q14").Value = 1666.431354929771
$ws.Range("R14").Value = 14997.88219436794
$ws.Range("S14").Value = 0.04774237628000955
$ws.Range("T14").Value = 0.04910356364522234

# Row 15: Inflammatory-Mac -> MuSCs
$ws.Range("G15").Value = 32.300192
$ws.Range("H15").Value = 96.900576
$ws.Range("I15").Value = 0.2476229975407503
$ws.Range("J15").Value = 0.2539410918713864
$ws.Range("M15").Value = 2.338518
$ws.Range("N15").Value = 4.677036
$ws.Range("O15").Value = 0.0087391953474509
$ws.Range("P15").Value = 0.005843151713055659
$ws.Range("Q15").Value = 75.534580395456
$ws.Range("R15").Value = 453.207482372736
$ws.Range("S15").Value = 0.00216402574802997
$ws.Range("T15").Value = 0.001483816325983516

# Row 16: Inflammatory-Mac -> Resolving-Mac
$ws.Range("G16").Value = 32.300192
$ws.Range("H16").Value = 96.900576
$ws.Range("I16").Value = 0.2476229975407503
$ws.Range("J16").Value = 0.2539410918713864
$ws.Range("M16").Value = 67.05398933333333
$ws.Range("N16").Value = 201.161968
$ws.Range("O16").Value = 0.2505851618888069
$ws.Range("P16").Value = 0.2513172654477853
$ws.Range("Q16").Value = 2165.856729832619
$ws.Range("R16").Value = 19492.71056849357
$ws.Range("S16").Value = 0.06205064892614053
$ws.Range("T16").Value = 0.06381978079394164

# Row 17: MuSCs -> ECs
$ws.Range("G17").Value = 9.736177999999999
$ws.Range("H17").Value = 19.472356
$ws.Range("I17").Value = 0.07464047213559308
$ws.Range("J17").Value = 0.0510299478916239
$ws.Range("M17").Value = 145.7087706666667
$ws.Range("N17").Value = 437.126312
$ws.Range("O17").Value = 0.5445232453600627
$ws.Range("P17").Value = 0.5461141113270247
$ws.Range("Q17").Value = 1418.646527371845
$ws.Range("R17").Value = 8511.879164231072
$ws.Range("S17").Value = 0.04064347212248047
$ws.Range("T17").Value = 0.02786817464389856

# Row 18: MuSCs -> FAPs
$ws.Range("G18").Value = 9.736177999999999
$ws.Range("H18").Value = 19.472356
$ws.Range("I18").Value = 0.07464047213559308
$ws.Range("J18").Value = 0.0510299478916239
$ws.Range("O18").Value = 0.003349722554576428
$ws.Range("P18").Value = 0.003359509023117945
$ws.Range("Q18").Value = 8.727032886478
$ws.Range("R18").Value = 52.36219731886799
$ws.Range("S18").Value = 0.0002500248729968296
$ws.Range("T18").Value = 0.000171435570391149

# Row 19: MuSCs -> Inflammatory-Mac
$ws.Range("G19").Value = 9.736177999999999
$ws.Range("H19").Value = 19.472356
$ws.Range("I19").Value = 0.07464047213559308
$ws.Range("J19").Value = 0.0510299478916239
$ws.Range("M19").Value = 51.59199533333333
$ws.Range("N19").Value = 154.775986
$ws.Range("O19").Value = 0.1928026748491032
$ws.Range("P19").Value = 0.1933659624890163
$ws.Range("Q19").Value = 502.3088499405026
$ws.Range("R19").Value = 3013.853099643015
$ws.Range("S19").Value = 0.0143908826797423
$ws.Range("T19").Value = 0.009867454989828202

# Row 20: MuSCs -> MuSCs
$ws.Range("G20").Value = 9.736177999999999
$ws.Range("H20").Value = 19.472356
$ws.Range("I20").Value = 0.07464047213559308
$ws.Range("J20").Value = 0.0510299478916239
$ws.Range("M20").Value = 2.338518
$ws.Range("N20").Value = 4.677036
$ws.Range("O20").Value = 0.0087391953474509
$ws.Range("P20").Value = 0.005843151713055659
$ws.Range("Q20").Value = 22.768227504204
$ws.Range("R20").Value = 91.072910016816
$ws.Range("S20").Value = 0.0006522976668189135
$ws.Range("T20").Value = 0.0002981757274400832

# Row 21: MuSCs -> Resolving-Mac
$ws.Range("G21").Value = 9.736177999999999
$ws.Range("H21").Value = 19.472356
$ws.Range("I21").Value = 0.07464047213559308
$ws.Range("J21").Value = 0.0510299478916239
$ws.Range("M21").Value = 67.05398933333333
$ws.Range("N21").Value = 201.161968
$ws.Range("O21").Value = 0.2505851618888069
$ws.Range("P21").Value = 0.2513172654477853
$ws.Range("Q21").Value = 652.8495757594346
$ws.Range("R21").Value = 3917.097454556608
$ws.Range("S21").Value = 0.01870379479355457
$ws.Range("T21").Value = 0.01282470696006589

# Row 22: Resolving-Mac -> ECs
$ws.Range("G22").Value = 36.79199966666667
$ws.Range("H22").Value = 110.375999
$ws.Range("I22").Value = 0.2820585476079611
$ws.Range("J22").Value = 0.2892552641013719
$ws.Range("M22").Value = 145.7087706666667
$ws.Range("N22").Value = 437.126312
$ws.Range("O22").Value = 0.5445232453600627
$ws.Range("P22").Value = 0.5461141113270247
$ws.Range("Q22").Value = 5360.91704179841
$ws.Range("R22").Value = 48248.25337618569
$ws.Range("S22").Value = 0.1535874357250327
$ws.Range("T22").Value = 0.1579663815013845

# Row 23: Resolving-Mac -> FAPs
$ws.Range("G23").Value = 36.79199966666667
$ws.Range("H23").Value = 110.375999
$ws.Range("I23").Value = 0.2820585476079611
$ws.Range("J23").Value = 0.2892552641013719
$ws.Range("O23").Value = 0.003349722554576428
$ws.Range("P23").Value = 0.003359509023117945
$ws.Range("Q23").Value = 32.97854569321633
$ws.Range("R23").Value = 296.806911238947
$ws.Range("S23").Value = 0.0009448178786334567
$ws.Range("T23").Value = 0.0009717556697329229

# Row 24: Resolving-Mac -> Inflammatory-Mac
$ws.Range("G24").Value = 36.79199966666667
$ws.Range("H24").Value = 110.375999
$ws.Range("I24").Value = 0.2820585476079611
$ws.Range("J24").Value = 0.2892552641013719
$ws.Range("M24").Value = 51.59199533333333
$ws.Range("N24").Value = 154.775986
$ws.Range("O24").Value = 0.1928026748491032
$ws.Range("P24").Value = 0.1933659624890163
$ws.Range("Q24").Value = 1898.172675106668
$ws.Range("R24").Value = 17083.55407596002
$ws.Range("S24").Value = 0.05438164244286802
$ws.Range("T24").Value = 0.05593212254797637

# Row 25: Resolving-Mac -> MuSCs
$ws.Range("G25").Value = 36.79199966666667
$ws.Range("H25").Value = 110.375999
$ws.Range("I25").Value = 0.2820585476079611
$ws.Range("J25").Value = 0.2892552641013719
$ws.Range("M25").Value = 2.338518
$ws.Range("N25").Value = 4.677036
$ws.Range("O25").Value = 0.0087391953474509
$ws.Range("P25").Value = 0.005843151713055659
$ws.Range("Q25").Value = 86.03875347649401
$ws.Range("R25").Value = 516.2325208589641
$ws.Range("S25").Value = 0.002464964746964252
$ws.Range("T25").Value = 0.001690162391944298

# Row 26: Resolving-Mac -> Resolving-Mac
$ws.Range("G26").Value = 36.79199966666667
$ws.Range("H26").Value = 110.375999
$ws.Range("I26").Value = 0.2820585476079611
$ws.Range("J26").Value = 0.2892552641013719
$ws.Range("M26").Value = 67.05398933333333
$ws.Range("N26").Value = 201.161968
$ws.Range("O26").Value = 0.2505851618888069
$ws.Range("P26").Value = 0.2513172654477853
$ws.Range("Q26").Value = 2467.050353200671
$ws.Range("R26").Value = 22203.45317880603
$ws.Range("S26").Value = 0.07067968681446267
$ws.Range("T26").Value = 0.0726948419903337

